$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (shared-string text updates only) ---
$ws.Range("C1").Value = "player_id"
$ws.Range("E1").Value = "birth_year"

# --- Re-seed the whole data block (rows 2-14, columns A:AM) with the
#     post-merge values. Column B (season_ending_year) is textual in the
#     original data (e.g. "1962"), so each value is apostrophe-prefixed to
#     keep Excel from coercing it to a number. ---
$data = @(
@(2332,"'1962",415,"Bill Russell",1934,"C",28,6,"NBA","BOS",76,$null,45.2,7.6,16.6,0.457,$null,$null,$null,7.6,16.6,0.457,0.457,3.8,6.3,0.595,$null,$null,23.6,4.5,"N/A - Stat tracked as of the 1973-74 NBA Season","N/A - Stat tracked as of the 1973-74 NBA Season",$null,2.7,18.9,"1960-62","Yes","No",1962),
@(2600,"'1964",415,"Bill Russell",1934,"C",30,8,"NBA","BOS",78,$null,44.6,6,13.8,0.433,$null,$null,$null,6,13.8,0.433,0.433,3,5.5,0.55,$null,$null,24.7,4.7,"N/A - Stat tracked as of the 1973-74 NBA Season","N/A - Stat tracked as of the 1973-74 NBA Season",$null,2.4,15,"1963-64","Yes","No",1964),
@(2467,"'1963",415,"Bill Russell",1934,"C",29,7,"NBA","BOS",78,$null,44.9,6.6,15.2,0.432,$null,$null,$null,6.6,15.2,0.432,0.432,3.7,6.6,0.555,$null,$null,23.6,4.5,"N/A - Stat tracked as of the 1973-74 NBA Season","N/A - Stat tracked as of the 1973-74 NBA Season",$null,2.4,16.8,"1962-63","Yes","No",1963),
@(2735,"'1965",415,"Bill Russell",1934,"C",31,9,"NBA","BOS",78,$null,44.4,5.5,12.6,0.438,$null,$null,$null,5.5,12.6,0.438,0.438,3.1,5.5,0.573,$null,$null,24.1,5.3,"N/A - Stat tracked as of the 1973-74 NBA Season","N/A - Stat tracked as of the 1973-74 NBA Season",$null,2.6,14.1,"1964-65","Yes","No",1965),
@(2995,"'1967",415,"Bill Russell",1934,"C",33,11,"NBA","BOS",81,$null,40.7,4.9,10.7,0.454,$null,$null,$null,4.9,10.7,0.454,0.454,3.5,5.8,0.61,$null,$null,21,5.8,"N/A - Stat tracked as of the 1973-74 NBA Season","N/A - Stat tracked as of the 1973-74 NBA Season",$null,3.2,13.3,"1966-67","Yes","No",1967),
@(1875,"'1958",415,"Bill Russell",1934,"C",24,2,"NBA","BOS",69,$null,38.3,6.6,15,0.442,$null,$null,$null,6.6,15,0.442,0.442,3.3,6.4,0.519,$null,$null,22.7,2.9,"N/A - Stat tracked as of the 1973-74 NBA Season","N/A - Stat tracked as of the 1973-74 NBA Season",$null,2.6,16.6,"1957-58","Yes","No",1958),
@(2219,"'1961",415,"Bill Russell",1934,"C",27,5,"NBA","BOS",78,$null,44.3,6.8,16,0.426,$null,$null,$null,6.8,16,0.426,0.426,3.3,6,0.55,$null,$null,23.9,3.4,"N/A - Stat tracked as of the 1973-74 NBA Season","N/A - Stat tracked as of the 1973-74 NBA Season",$null,2,16.9,"1960-61","Yes","No",1961),
@(2865,"'1966",415,"Bill Russell",1934,"C",32,10,"NBA","BOS",78,$null,43.4,5,12.1,0.415,$null,$null,$null,5,12.1,0.415,0.415,2.9,5.2,0.551,$null,$null,22.8,4.8,"N/A - Stat tracked as of the 1973-74 NBA Season","N/A - Stat tracked as of the 1973-74 NBA Season",$null,2.8,12.9,"1965-66","Yes","No",1966),
@(3150,"'1968",415,"Bill Russell",1934,"C",34,12,"NBA","BOS",78,$null,37.9,4.7,11,0.425,$null,$null,$null,4.7,11,0.425,0.425,3.2,5.9,0.537,$null,$null,18.6,4.6,"N/A - Stat tracked as of the 1973-74 NBA Season","N/A - Stat tracked as of the 1973-74 NBA Season",$null,3.1,12.5,"1967-68","Yes","No",1968),
@(1760,"'1957",415,"Bill Russell",1934,"C",23,1,"NBA","BOS",48,$null,35.3,5.8,13.5,0.427,$null,$null,$null,5.8,13.5,0.427,0.427,3.2,6.4,0.492,$null,$null,19.6,1.8,"N/A - Stat tracked as of the 1973-74 NBA Season","N/A - Stat tracked as of the 1973-74 NBA Season",$null,3,14.7,"1956-57","Yes","No",1957),
@(2097,"'1960",415,"Bill Russell",1934,"C",26,4,"NBA","BOS",74,$null,42.5,7.5,16.1,0.467,$null,$null,$null,7.5,16.1,0.467,0.467,3.2,5.3,0.612,$null,$null,24,3.7,"N/A - Stat tracked as of the 1973-74 NBA Season","N/A - Stat tracked as of the 1973-74 NBA Season",$null,2.8,18.2,"1959-60","Yes","No",1960),
@(3513,"'1969",415,"Bill Russell",1934,"C",35,13,"NBA","BOS",77,$null,42.7,3.6,8.4,0.433,$null,$null,$null,3.6,8.4,0.433,0.433,2.6,5,0.526,$null,$null,19.3,4.9,"N/A - Stat tracked as of the 1973-74 NBA Season","N/A - Stat tracked as of the 1973-74 NBA Season",$null,3,9.9,"1968-69","No","No",1969),
@(1985,"'1959",415,"Bill Russell",1934,"C",25,3,"NBA","BOS",70,$null,42.6,6.5,14.2,0.457,$null,$null,$null,6.5,14.2,0.457,0.457,3.7,6.1,0.598,$null,$null,23,3.2,"N/A - Stat tracked as of the 1973-74 NBA Season","N/A - Stat tracked as of the 1973-74 NBA Season",$null,2.3,16.7,"1958-59","Yes","No",1959)
)

$arr = New-Object 'object[,]' 13,39
for ($r = 0; $r -lt 13; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt 39; $c++) {
        $arr[$r, $c] = $row[$c]
    }
}

$ws.Range("A2:AM14").Value = $arr

Write-Host "Applied Bill Russell per-game averages refresh"
